# Insert a new weekly record (date 2021-10-28 / serial 44497) as a new row 23,
# pushing the existing rows 23-104 down to 24-105 (dimension grows A1:R104 -> A1:R105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 23; this shifts rows 23..104 down
# to 24..105 and carries the row-above formatting (matches the D column's
# date style s="2" already present on the neighbouring rows).
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with the new weekly observation.
$ws.Cells.Item(23, 1).Value2  = 2
$ws.Cells.Item(23, 2).Value2  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(23, 3).Value2  = "Coquimbo"
$ws.Cells.Item(23, 4).Value2  = 44497
$ws.Cells.Item(23, 5).Value2  = 4
$ws.Cells.Item(23, 6).Value2  = 100112031
$ws.Cells.Item(23, 7).Value2  = "Poroto verde"
$ws.Cells.Item(23, 8).Value2  = "Magnum"
$ws.Cells.Item(23, 9).Value2  = "Primera"
$ws.Cells.Item(23, 10).Value2 = 1000
$ws.Cells.Item(23, 11).Value2 = 35000
$ws.Cells.Item(23, 12).Value2 = 40000
$ws.Cells.Item(23, 13).Value2 = 37500
$ws.Cells.Item(23, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(23, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(23, 16).Value2 = 1500
$ws.Cells.Item(23, 17).Value2 = 25
$ws.Cells.Item(23, 18).Value2 = "Hortaliza"
